$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").NumberFormat = "@"
$wsOverview.Range("G2").Value = "2016-09-02 20:19:00"
$wsOverview.Range("G5").NumberFormat = "@"
$wsOverview.Range("G5").Value = "2016-09-02 20:19:00"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").NumberFormat = "@"
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").NumberFormat = "@"
$wsZhCn.Range("E5").Value = "mt"

$wsZhCn.Range("H2").NumberFormat = "@"
$wsZhCn.Range("H2").Value = "2016-09-02 20:18:55"
$wsZhCn.Range("H5").NumberFormat = "@"
$wsZhCn.Range("H5").Value = "2016-09-02 20:18:55"

$wsZhCn.Range("K2").NumberFormat = "@"
$wsZhCn.Range("K2").Value = "2016-09-02 20:19:26"
$wsZhCn.Range("K5").NumberFormat = "@"
$wsZhCn.Range("K5").Value = "2016-09-02 20:19:26"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").NumberFormat = "@"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").NumberFormat = "@"
$wsDeDe.Range("E5").Value = "mt"

$wsDeDe.Range("H2").NumberFormat = "@"
$wsDeDe.Range("H2").Value = "2016-09-02 20:19:00"
$wsDeDe.Range("H5").NumberFormat = "@"
$wsDeDe.Range("H5").Value = "2016-09-02 20:19:00"

$wsDeDe.Range("K2").NumberFormat = "@"
$wsDeDe.Range("K2").Value = "2016-09-02 20:19:33"
